# Auto-generated edit script: refresh the cryptos price/volume table.
# For numeric-looking text values (e.g. "1.00", "215.43") a leading "'" forces
# Excel to keep them as text instead of coercing to a Double, and the Style
# reset afterwards drops the transient quote-prefix format so cell styling
# stays identical to the original (default) style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.162.20"
$ws.Range("E2").Value = "  +3.61%  "
$ws.Range("D3").Value = "1.659.80"
$ws.Range("E3").Value = "  +4.30%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'215.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.74%  "
$ws.Range("E6").Value = "  +1.08%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +2.20%  "
$ws.Range("D9").Value = "'0.0614"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.65%  "
$ws.Range("D10").Value = "'19.59"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.84%  "
$ws.Range("D11").Value = "'0.0863"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.23%  "
$ws.Range("D12").Value = "1.896.00"
$ws.Range("E12").Value = "  +4.48%  "
$ws.Range("D13").Value = "1.662.79"
$ws.Range("E13").Value = "  +4.65%  "
$ws.Range("D14").Value = "'4.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.80%  "
$ws.Range("D15").Value = "'0.520"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.40%  "
$ws.Range("D16").Value = "'64.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.14%  "
$ws.Range("D17").Value = "27.176.51"
$ws.Range("E17").Value = "  +3.77%  "
$ws.Range("D18").Value = "'238.85"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.66%  "
$ws.Range("D19").Value = "'7.85"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.03%  "
$ws.Range("D20").Value = "0.0₃0729"
$ws.Range("E20").Value = "  +1.44%  "
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("E22").Value = "  +4.99%  "
$ws.Range("E23").Value = "  +5.65%  "
$ws.Range("D24").Value = "'9.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.74%  "
$ws.Range("D25").Value = "'145.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.28%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "'7.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.39%  "
$ws.Range("E28").Value = "  +1.42%  "
$ws.Range("D29").Value = "'15.83"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.43%  "
$ws.Range("D30").Value = "'0.0497"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.11%  "
$ws.Range("E31").Value = "  +1.38%  "
$ws.Range("D32").Value = "1.539.01"
$ws.Range("E32").Value = "  +5.72%  "
$ws.Range("D33").Value = "'3.28"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.91%  "
$ws.Range("D34").Value = "'3.04"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.58%  "
$ws.Range("E35").Value = "  +8.25%  "
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("D37").Value = "'0.573"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.76%  "
$ws.Range("D38").Value = "'0.888"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.75%  "
$ws.Range("E39").Value = "  +3.27%  "
$ws.Range("D40").Value = "'5.95"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.52%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("E42").Value = "  +4.96%  "
$ws.Range("D43").Value = "'66.23"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +10.02%  "
$ws.Range("D44").Value = "1.802.28"
$ws.Range("E44").Value = "  +4.33%  "
$ws.Range("E45").Value = "  +2.80%  "
$ws.Range("D46").Value = "'0.921"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("D47").Value = "'89.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.83%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0105"
$ws.Range("E48").Value = "  +0.87%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'1.53"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.47%  "
$ws.Range("E50").Value = "  +0.98%  "
$ws.Range("D51").Value = "'0.0975"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.43%  "
